# Insert a new first worksheet "generate-weights-more-distance" ahead of the
# existing four sheets (generate-weights_test, datanoise-rate, weigths-sample,
# generate-weights). The new sheet becomes the active/selected tab (activeTab=0),
# matching the target workbook layout.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() with no arguments inserts the new sheet before the first
# existing sheet and makes it the active sheet (tabSelected=true, activeTab=0).
$ws = $wb.Worksheets.Add()
$ws.Name = "generate-weights-more-distance"

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "WeightName"
$ws.Range("B1").Value = "caseID"
$ws.Range("C1").Value = "weightVar"

# --- Named weight rows -------------------------------------------------
$ws.Range("A2").Value = "Ext-weight"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.05

$ws.Range("A3").Value = "Pos-weight"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 0.1

$ws.Range("A4").Value = "Neg-weight"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 0.2

# --- Remaining caseID / weightVar pairs (no WeightName label) ---------
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 0.3

$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 0.4

$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 0.5

$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 0.6

$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 0.7

$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 0.8

$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 0.9

$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 1

# --- Reserved (currently empty) column D, wrapped like the D column of
# the "generate-weights" sheet it mirrors -------------------------------
$ws.Range("D2:D6").WrapText = $true

# --- Column widths (match the custom widths on D/E from the source sheet)
$ws.Columns.Item(4).ColumnWidth = 13.47
$ws.Columns.Item(5).ColumnWidth = 12.37

# --- Selection matches the saved cursor position in the target sheet ---
$ws.Range("D4").Select()
